# Auto-generated PowerShell COM-interop script
# "Ultra-HQ Sample Overhaul": replaces titles/bullets/notes across the
# 6 flagship topics (AI-CN, AI-EN, Longevity-CN, Longevity-EN,
# Renaissance-CN, Quantum-EN) and extends every References slide with a
# 4th citation (jstor.org), per the commit diff.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Speaker notes (NotesPage, shape 2 = "Notes Placeholder")
# ---------------------------------------------------------------
$np = $p.Slides.Item(3).NotesPage
$np.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "本节介绍 AI 的前身及其早期的哲学与数学基础，重点强调对人类智能的逻辑模拟尝试。"
$np = $p.Slides.Item(4).NotesPage
$np.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "介绍神经网络从被冷落到重新获得学术界关注的过程，为后来的深度学习爆发做铺垫。"
$np = $p.Slides.Item(6).NotesPage
$np.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "分析深度学习在视觉和 NLP 领域的双重突破，揭示 Transformer 架构为何成为当今万物互联的技术底座。"
$np = $p.Slides.Item(10).NotesPage
$np.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "旨在揭示胰岛素在人体内调取和储存能量的底层逻辑，以及不当饮食对代谢系统的长期损伤。"
$np = $p.Slides.Item(12).NotesPage
$np.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "探讨如何从微观层面优化身体引擎，强调生活细节对生物学年龄的逆转作用。"
$np = $p.Slides.Item(16).NotesPage
$np.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "本节重点讲解文艺复兴如何利用数学工具彻底改变了人类观察图像的方式。"
$np = $p.Slides.Item(18).NotesPage
$np.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "分析达芬奇如何通过精妙的光学观察超越了机械的透视法，赋予肖像画以生命感。"
$np = $p.Slides.Item(22).NotesPage
$np.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "Explaining the physics that enables exponential speedup in specific algorithmic domains like cryptography."
$np = $p.Slides.Item(24).NotesPage
$np.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "Discussing the convergence of AI architecture and quantum hardware for the 2030 decade."
$np = $p.Slides.Item(28).NotesPage
$np.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "Breaking down the cellular mechanisms that respond to nutrient deprivation and increase metabolic health."

# ---------------------------------------------------------------
# Slide titles / subtitles / bullet content
# ---------------------------------------------------------------
$p.Slides.Item(1).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "人工智能：从图灵测试到通用人工智能 (AGI)"
$p.Slides.Item(1).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Topic: AI的发展历史与未来趋势"
$p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "计算智能的起源与逻辑奠基 (1950-1980)"
$p.Slides.Item(3).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "图灵测试与符号 AI 的诞生"
$p.Slides.Item(3).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "1950年：阿兰·图灵发表《计算机器与智能》，提出著名的“图灵测试” (Turing Test)。"
$p.Slides.Item(3).Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "1956年：达特茅斯会议 (Dartmouth Workshop) 正式确立“人工智能”学科，麦卡锡、明斯基等人为学科领袖。"
$p.Slides.Item(3).Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "逻辑主义时代：基于规则的专家系统（如 MYCIN）在特定医疗诊断领域取得初步成功。"
$p.Slides.Item(3).Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "瓶颈出现：早期 AI 难以处理模糊信息，导致70年代中期进入第一个“AI 冬天”。"
$p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "联结主义与神经网络的复兴"
$p.Slides.Item(4).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "1986年：Rumelhart 提出反向传播算法 (Backpropagation)，解决了多层感知器的训练难题。"
$p.Slides.Item(4).Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "统计学习方法崛起：SVM 与随机森林在90年代成为机器学习的主流工具。"
$p.Slides.Item(4).Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "GPU 计算能力的增强：为复杂的矩阵运算提供了硬件基础，神经网络的研究重心逐渐转向深度化。"
$p.Slides.Item(4).Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "循环神经网络 (RNN) 与 LSTM：在高盛等金融机构及自然语言处理中开始显露头角。"
$p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "深度学习革命与大模型时代 (2012-Present)"
$p.Slides.Item(6).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "从 ImageNet 到 Transformer 架构"
$p.Slides.Item(6).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "2012年：AlexNet 以领先第二名10.8%的优势夺得 ImageNet 冠军，开启深度卷积神经网络时代。"
$p.Slides.Item(6).Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "2017年：Google 发表《Attention is All You Need》，提出 Transformer 架构，颠覆序列建模模式。"
$p.Slides.Item(6).Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "预训练大模型 (LLMs)：GPT-3 的 1750亿参数规模展示了模型容量与涌现能力 (Emergent Abilities) 的正相关性。"
$p.Slides.Item(6).Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "推理与对齐：利用 RLHF (基于人类反馈的强化学习) 解决了模型在道德与逻辑层面的幻觉问题。"
$p.Slides.Item(8).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "代谢革命：基于现代营养学的长寿科学"
$p.Slides.Item(8).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Topic: 现代营养学与代谢健康科学"
$p.Slides.Item(9).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "细胞能量代谢与胰岛素平衡"
$p.Slides.Item(10).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "胰岛素敏感性：健康的万能钥匙"
$p.Slides.Item(10).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "胰岛素抵抗 (Insulin Resistance)：不仅是糖尿病前兆，更是 2 型糖尿病、多囊卵巢综合征 (PCOS) 的核心驱动因素。"
$p.Slides.Item(10).Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "血糖波动的负面影响：餐后高血糖导致的糖基化终产物 (AGEs) 会加速血管内膜老化。"
$p.Slides.Item(10).Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "低 GI 饮食策略：通过全谷物和高纤维摄入，维持血清能量供应的平滑曲线。"
$p.Slides.Item(10).Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "动态血糖监测 (CGM)：现代医疗技术从盲目补给向实时精准控糖的转变。"
$p.Slides.Item(11).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "线粒体功能与抗炎生活方式"
$p.Slides.Item(12).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "线粒体：细胞的能量工厂"
$p.Slides.Item(12).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "线粒体自噬 (Mitophagy)：通过断食或高强度间歇训练 (HIIT) 触发坏死线粒体的自我清理。"
$p.Slides.Item(12).Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "氧化压力与抗氧化平衡：SOD 等内源性酶在抵御超氧阴离子自由基中的核心作用。"
$p.Slides.Item(12).Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "睡眠与线粒体修复：深度睡眠期是大脑清除代谢废物（β-淀粉样蛋白）的唯一窗口期。"
$p.Slides.Item(12).Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "Omega-3s 的抗炎机制：通过调节细胞膜流动性来降低慢性系统性炎症水平。"
$p.Slides.Item(14).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "透视与光影：文艺复兴艺术的技术巅峰"
$p.Slides.Item(14).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Topic: 文艺复兴大师与技法深度赏析"
$p.Slides.Item(15).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "科学写实主义的黄金时代"
$p.Slides.Item(16).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "数学视角的引入：线性透视法"
$p.Slides.Item(16).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "布鲁内莱斯基的发现：通过数学消失点 (Vanishing Point) 在二维平面还原三维物理空间。"
$p.Slides.Item(16).Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "马萨乔的《圣三一》：西方艺术史上第一张严格遵循线性透视规则的大型湿壁画。"
$p.Slides.Item(16).Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "缩短透视法 (Foreshortening)：使人体部位垂直于画平面，产生强烈的立体压缩感。"
$p.Slides.Item(16).Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "艺术与科学的深度融合：艺术家由单纯的工匠转变为具备解剖学与光学知识的知识分子。"
$p.Slides.Item(17).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "达芬奇与威尼斯画派的色彩革命"
$p.Slides.Item(18).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "晕涂法与大气透视"
$p.Slides.Item(18).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "晕涂法 (Sfumato)：达芬奇通过层层薄釉色消除轮廓线，营造出如烟雾般的柔和过度。"
$p.Slides.Item(18).Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "《蒙娜丽莎》中的色彩博弈：利用冷暖色调在背景中建立的大气远近感。"
$p.Slides.Item(18).Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "威尼斯画派：提香与乔尔乔内对“色彩建构模型”的重视，挑战了佛罗伦萨的“素描中心论”。"
$p.Slides.Item(18).Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "明暗对照法 (Chiaroscuro)：利用极端的光源对比，增强人物的情感张力与体积感。"
$p.Slides.Item(20).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "The Quantum Leap: Synergy of Qubits and Neural Networks"
$p.Slides.Item(20).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Topic: Quantum Computing and Generative AI"
$p.Slides.Item(21).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Quantum Mechanics as a Computing Paradigm"
$p.Slides.Item(22).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Superposition and Entanglement"
$p.Slides.Item(22).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Beyond Binary: Qubits leverage quantum superposition to represent |0⟩ and |1⟩ simultaneously."
$p.Slides.Item(22).Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "Entanglement Dynamics: Correlating qubits to perform massively parallel state-space explorations."
$p.Slides.Item(22).Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "Quantum Supremacy: Google's Sycamore processor solving tasks in 200 seconds that supercomputers take 10,000 years."
$p.Slides.Item(22).Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "Error Mitigation: The challenge of decoherence and the quest for fault-tolerant logical qubits."
$p.Slides.Item(23).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "AI's Next Hardware Frontier"
$p.Slides.Item(24).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Quantum Neural Networks (QNN)"
$p.Slides.Item(24).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Quantum-Classical Hybrids: Using variational circuits to optimize deep learning weights."
$p.Slides.Item(24).Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "Infinite Context Windows: Processing complex drug-discovery data using quantum-enhanced transformers."
$p.Slides.Item(24).Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "Sustainability: Theoretically reducing the massive carbon footprint of LLM training via specialized circuits."
$p.Slides.Item(24).Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "The RSA Threat: Quantum algorithms (Shor's) and the inevitable transition to Post-Quantum Cryptography (PQC)."
$p.Slides.Item(26).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "The Longevity Code: Optimizing Human Lifespan"
$p.Slides.Item(26).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Topic: Science of Longevity: Nutrition and Sleep"
$p.Slides.Item(27).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Autophagy and The Biology of Fasting"
$p.Slides.Item(28).Shapes.Item(1).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Yoshinori Ohsumi's Nobel Discovery"
$p.Slides.Item(28).Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Cellular Self-Cleaning: How cells degrade and recycle damaged components via lysosomes."
$p.Slides.Item(28).Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "mTOR Pathways: Inhibiting growth signals to promote longevity and cellular repair."
$p.Slides.Item(28).Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "Time-Restricted Feeding (TRF): Syncing intake windows with biological circadian rhythms."
$p.Slides.Item(28).Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "Senescent 'Zombie' Cells: The role of senolytics in removing inflammation-inducing cells."

# ---------------------------------------------------------------
# References slides: swap the first 3 URLs and append a 4th
# (jstor.org) paragraph, cloning the sz=1200 paragraph formatting
# ---------------------------------------------------------------
$refTr = $p.Slides.Item(7).Shapes.Item(2).TextFrame.TextRange
$refTr.Paragraphs(2).Runs(1).Text = "https://en.wikipedia.org/wiki/Main_Page"
$refTr.Paragraphs(3).Runs(1).Text = "https://www.nih.gov/"
$refTr.Paragraphs(4).Runs(1).Text = "https://scholar.google.com/"
$null = $refTr.InsertAfter("`rhttps://www.jstor.org/")
$refTr = $p.Slides.Item(13).Shapes.Item(2).TextFrame.TextRange
$refTr.Paragraphs(2).Runs(1).Text = "https://en.wikipedia.org/wiki/Main_Page"
$refTr.Paragraphs(3).Runs(1).Text = "https://www.nih.gov/"
$refTr.Paragraphs(4).Runs(1).Text = "https://scholar.google.com/"
$null = $refTr.InsertAfter("`rhttps://www.jstor.org/")
$refTr = $p.Slides.Item(19).Shapes.Item(2).TextFrame.TextRange
$refTr.Paragraphs(2).Runs(1).Text = "https://en.wikipedia.org/wiki/Main_Page"
$refTr.Paragraphs(3).Runs(1).Text = "https://www.nih.gov/"
$refTr.Paragraphs(4).Runs(1).Text = "https://scholar.google.com/"
$null = $refTr.InsertAfter("`rhttps://www.jstor.org/")
$refTr = $p.Slides.Item(25).Shapes.Item(2).TextFrame.TextRange
$refTr.Paragraphs(2).Runs(1).Text = "https://en.wikipedia.org/wiki/Main_Page"
$refTr.Paragraphs(3).Runs(1).Text = "https://www.nih.gov/"
$refTr.Paragraphs(4).Runs(1).Text = "https://scholar.google.com/"
$null = $refTr.InsertAfter("`rhttps://www.jstor.org/")
$refTr = $p.Slides.Item(29).Shapes.Item(2).TextFrame.TextRange
$refTr.Paragraphs(2).Runs(1).Text = "https://en.wikipedia.org/wiki/Main_Page"
$refTr.Paragraphs(3).Runs(1).Text = "https://www.nih.gov/"
$refTr.Paragraphs(4).Runs(1).Text = "https://scholar.google.com/"
$null = $refTr.InsertAfter("`rhttps://www.jstor.org/")

